# Applies the "Saldo" export refresh:
#  - rewrites the top block of the Export sheet (rows 2-10) with the
#    refreshed data (re-ordered, some balances updated)
#  - removes the now-superseded rows that duplicated the new block
#    (old rows 11-15), since the sheet shrank by 5 rows overall
#
# Resulting top of sheet (A:Conta, B:Nome, C:Saldo):
#   1  Conta       Nome       Saldo   (header, unchanged)
#   2  005064129   THIAGO     20349.02
#   3  004461070   EDUARDO    16152.95
#   4  004392159   RODRIGO    900.21
#   5  005696595   CLUBE      785.26
#   6  004574428   GUILHERME  745.08
#   7  004488571   CARLOS     440.36
#   8  004322719   GISELA     276.97
#   9  004382374   THEOMAR    100.74
#  10  004459875   HELVECIO   100.57
#  11  002823185   SIMONE     100.22   (already there - untouched, just shifts up)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# New content for rows 2 through 10 (Conta / Nome / Saldo)
$newRows = @(
    @("005064129", "THIAGO",    20349.02),
    @("004461070", "EDUARDO",   16152.95),
    @("004392159", "RODRIGO",   900.21),
    @("005696595", "CLUBE",     785.26),
    @("004574428", "GUILHERME", 745.08),
    @("004488571", "CARLOS",    440.36),
    @("004322719", "GISELA",    276.97),
    @("004382374", "THEOMAR",   100.74),
    @("004459875", "HELVECIO",  100.57)
)

# Make sure the account-number column keeps its leading zeros (it is
# plain "General" formatted in the source file, so a bare numeric-looking
# string would otherwise be auto-coerced to a number on assignment).
$ws.Range("A2:A10").NumberFormat = "@"

$r = 2
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Old rows 11-15 (RODRIGO, GUILHERME, CARLOS, GISELA, THEOMAR) now duplicate
# the rows we just wrote above, so delete that whole block; everything below
# (SIMONE, ...) shifts up to close the gap.
$ws.Range("A11:C15").EntireRow.Delete()
